# Fix minor bugs and add new layout sizes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 (the single data row) with the corrected / new design values
$ws.Range("A2").Value = 13.668482276188673   # weight
$ws.Range("B2").Value = 0.94576720880508192  # S_w
$ws.Range("C2").Value = 2                    # b_w
$ws.Range("D2").Value = 4.2293705710665854   # A
$ws.Range("J2").Value = 3                    # L_fuse
$ws.Range("M2").Value = 0.15                 # S_ht
$ws.Range("O2").Value = 0.5                  # b_h
$ws.Range("Q2").Value = 0.4                  # S_vt
$ws.Range("V2").Value = 14.452269172514399   # W_S
$ws.Range("W2").Value = 0.13367625320445395  # Preq_W
$ws.Range("X2").Value = 1.8271514976723882   # P_needed
